# C5-PowerPoint.pptx edit
#
# 1) The table on slide 6 ("SOURCES OF FINANCE") switches from the custom
#    table style {363A3CA6-E62E-4840-AC3F-6E83845CCF75} to the built-in
#    table style {DC1A7209-A95B-41C4-9B40-11BB06270EA6}.
#
# 2) The deck's theme is switched from the "Integral" design back to the
#    default "Office Theme" design (the Office Theme colors are written
#    into the active theme's color scheme through the exposed
#    ThemeColorScheme object model).

$p = $ppt.ActivePresentation

# --- 1) Re-style the table on slide 6 -------------------------------------
$slide = $p.Slides.Item(6)
$tableShape = $slide.Shapes.Item(2)
$table = $tableShape.Table
$table.ApplyStyle("{DC1A7209-A95B-41C4-9B40-11BB06270EA6}")

# --- 2) Switch the design's colour scheme back to "Office Theme" ---------
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink (this is the order used by
# ThemeColorScheme.Colors / MsoThemeColorSchemeIndex).
$officeThemeColors = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000",  # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

$themeColors = $slide.ThemeColorScheme
for ($i = 1; $i -le $officeThemeColors.Length; $i++) {
    $hex = $officeThemeColors[$i - 1]
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    # OLE COLORREF long storage is 0x00BBGGRR.
    $themeColors.Colors($i).RGB = ($b * 65536) + ($g * 256) + $r
}
